# The NET-001 "Segmentation of Endpoint and Medical Devices" control row
# (row 4) was removed from the controls table entirely; every row below it
# shifts up by one. Deleting the whole row (not just clearing cell
# contents) removes the row from the sheet, re-numbers the remaining rows,
# and drops the three shared strings that were only referenced by that row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(4).Delete()

# Match the saved view state: zoomed out to 49% with cell J18 selected
# (and no frozen/scrolled topLeftCell override).
$excel.ActiveWindow.Zoom = 49
$ws.Range("J18").Select()
